$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 816.6667
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H8").Value = 86.666664
$ws.Range("I8").Value = 86.666664
$ws.Range("K8").Value = 259.999992
$ws.Range("M8").Value = -120.999992

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

$ws.Range("H100").Value = 4400
$ws.Range("I100").Value = 5033.3335
$ws.Range("K100").Value = 5033.3335
$ws.Range("M100").Value = -4492.3335

$ws.Range("H111").Value = 1294.8889
$ws.Range("I111").Value = 1044.25
$ws.Range("J111").Value = 3300
$ws.Range("K111").Value = 3132.75
$ws.Range("L111").Value = 9900
$ws.Range("M111").Value = -65.75
$ws.Range("N111").Value = -16034

$ws.Range("H132").Value = 4668.5
$ws.Range("I132").Value = 1232.5714
$ws.Range("K132").Value = 3697.7142
$ws.Range("M132").Value = -1167.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 20833.334
$ws.Range("I31").Value = 20833.334
$ws.Range("K31").Value = 20833.334
$ws.Range("M31").Value = -20539.334

$ws.Range("H45").Value = 3999.75
$ws.Range("I45").Value = 1499.5
$ws.Range("K45").Value = 1499.5
$ws.Range("M45").Value = -1122.5

$ws.Range("H61").Value = 8013
$ws.Range("I61").Value = 8182
$ws.Range("K61").Value = 8182
$ws.Range("M61").Value = -7970

$ws.Range("H74").Value = 609.6111
$ws.Range("I74").Value = 609.6111
$ws.Range("K74").Value = 609.6111
$ws.Range("M74").Value = 264.3889

$ws.Range("H77").Value = 609.6111
$ws.Range("I77").Value = 609.6111
$ws.Range("K77").Value = 3048.0555
$ws.Range("M77").Value = 1319.9445

$ws.Range("H97").Value = 10113.167
$ws.Range("I97").Value = 142.25
$ws.Range("K97").Value = 142.25
$ws.Range("M97").Value = 353.75

$ws.Range("H110").Value = 1667.1111
$ws.Range("I110").Value = 1625.5
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 1625.5
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 419.5
$ws.Range("N110").Value = -6090

$ws.Range("H136").Value = 8013
$ws.Range("I136").Value = 8182
$ws.Range("K136").Value = 24546
$ws.Range("M136").Value = -21996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 616.6111
$ws.Range("I94").Value = 640.5294
$ws.Range("K94").Value = 640.5294
$ws.Range("M94").Value = -189.5294

$ws.Range("H102").Value = 12649
$ws.Range("I102").Value = 12649
$ws.Range("K102").Value = 12649
$ws.Range("M102").Value = -9404

$ws.Range("H105").Value = 5155.4546
$ws.Range("I105").Value = 4492.2
$ws.Range("K105").Value = 4492.2
$ws.Range("M105").Value = -2745.2

$ws.Range("H107").Value = 1049.4117
$ws.Range("I107").Value = 1092.875
$ws.Range("J107").Value = 354
$ws.Range("K107").Value = 1092.875
$ws.Range("L107").Value = 354
$ws.Range("M107").Value = 827.125
$ws.Range("N107").Value = -4194

$ws.Range("H134").Value = 2465.353
$ws.Range("I134").Value = 1707.9286
$ws.Range("K134").Value = 5123.7858
$ws.Range("M134").Value = -2588.7858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1153.3334
$ws.Range("I7").Value = 1026.7
$ws.Range("J7").Value = 1311.625
$ws.Range("K7").Value = 1026.7
$ws.Range("L7").Value = 1311.625
$ws.Range("M7").Value = -913.7
$ws.Range("N7").Value = -1537.625

$ws.Range("H105").Value = 1997.1666
$ws.Range("I105").Value = 1494.3334
$ws.Range("K105").Value = 1494.3334
$ws.Range("M105").Value = 252.6666

$ws.Range("H141").Value = 137499.25
$ws.Range("J141").Value = 137499.25
$ws.Range("L141").Value = 137499.25
$ws.Range("N141").Value = -147859.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1271.5
$ws.Range("I5").Value = 1215.5
$ws.Range("K5").Value = 3646.5
$ws.Range("M5").Value = -3534.5

$ws.Range("H60").Value = 2324.25
$ws.Range("I60").Value = 149.5
$ws.Range("J60").Value = 4499
$ws.Range("K60").Value = 448.5
$ws.Range("L60").Value = 13497
$ws.Range("M60").Value = -197.5
$ws.Range("N60").Value = -13999

$ws.Range("H81").Value = 8400
$ws.Range("J81").Value = 12500
$ws.Range("L81").Value = 37500
$ws.Range("N81").Value = -39746

$ws.Range("H84").Value = 8400
$ws.Range("J84").Value = 12500
$ws.Range("L84").Value = 112500
$ws.Range("N84").Value = -123732

$ws.Range("H135").Value = 1271.5
$ws.Range("I135").Value = 1215.5
$ws.Range("K135").Value = 10939.5
$ws.Range("M135").Value = -8404.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 134999.5
$ws.Range("I10").Value = 19999
$ws.Range("J10").Value = 250000
$ws.Range("K10").Value = 19999
$ws.Range("L10").Value = 250000
$ws.Range("M10").Value = -19830
$ws.Range("N10").Value = -250338

$ws.Range("H22").Value = 3371
$ws.Range("I22").Value = 2556.5
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 2556.5
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -2027.5
$ws.Range("N22").Value = -6058

$ws.Range("H97").Value = 1029.909
$ws.Range("I97").Value = 1054.375
$ws.Range("J97").Value = 964.6667
$ws.Range("K97").Value = 1054.375
$ws.Range("L97").Value = 964.6667
$ws.Range("M97").Value = -558.375
$ws.Range("N97").Value = -1956.6667

$ws.Range("H107").Value = 794.8
$ws.Range("I107").Value = 531
$ws.Range("J107").Value = 1850
$ws.Range("K107").Value = 531
$ws.Range("L107").Value = 1850
$ws.Range("M107").Value = 1389
$ws.Range("N107").Value = -5690

$ws.Range("H122").Value = 2548.875
$ws.Range("I122").Value = 2899.5
$ws.Range("J122").Value = 1497
$ws.Range("K122").Value = 8698.5
$ws.Range("L122").Value = 4491
$ws.Range("M122").Value = -6248.5
$ws.Range("N122").Value = -9391

$ws.Range("H132").Value = 1777.9395
$ws.Range("I132").Value = 1778.931
$ws.Range("J132").Value = 1770.75
$ws.Range("K132").Value = 5336.793
$ws.Range("L132").Value = 5312.25
$ws.Range("M132").Value = -2806.793
$ws.Range("N132").Value = -10372.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws.Range("H17").Value = 22000
$ws.Range("J17").Value = 22000
$ws.Range("L17").Value = 22000
$ws.Range("N17").Value = -22340

$ws.Range("H25").Value = 19000
$ws.Range("I25").Value = 19000
$ws.Range("K25").Value = 19000
$ws.Range("M25").Value = -18770

$ws.Range("H46").Value = 1452.909
$ws.Range("I46").Value = 941.17645
$ws.Range("J46").Value = 1996.625
$ws.Range("K46").Value = 941.17645
$ws.Range("L46").Value = 1996.625
$ws.Range("M46").Value = -753.17645
$ws.Range("N46").Value = -2372.625

$ws.Range("H74").Value = 47499.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 47499.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 47499.5
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -49495.5

$ws.Range("H77").Value = 47499.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 47499.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 142498.5
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -152482.5

$ws.Range("H132").Value = 3423
$ws.Range("I132").Value = 3423
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10269
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7739
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 4000000
$ws.Range("I11").Value = 4000000
$ws.Range("K11").Value = 4000000
$ws.Range("M11").Value = -3999858

$ws.Range("H126").Value = 1319.7778
$ws.Range("I126").Value = 1297.25
$ws.Range("K126").Value = 3891.75
$ws.Range("M126").Value = -1421.75

$ws.Range("H127").Value = 24000
$ws.Range("I127").Value = 24000
$ws.Range("K127").Value = 24000
$ws.Range("M127").Value = -19040
